$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'basketball with design'
$ws.Cells.Item(2, 1).Value = 'yoga knee'
$ws.Cells.Item(3, 1).Value = 'basketball skins'
$ws.Cells.Item(4, 1).Value = 'leg pads hockey'
$ws.Cells.Item(5, 1).Value = 'elastic calf sleeve'
$ws.Cells.Item(6, 1).Value = 'softball shorts girls'
$ws.Cells.Item(7, 1).Value = 'compression knee sleeve with pads'
$ws.Cells.Item(8, 1).Value = 'volleyball compression shorts'
$ws.Cells.Item(9, 1).Value = 'knee compression sleeve for basketball'
$ws.Cells.Item(10, 1).Value = 'youth compression sleeve baseball'
$ws.Cells.Item(11, 1).Value = 'youth basketball shorts'
$ws.Cells.Item(12, 1).Value = 'hockey pants'
$ws.Cells.Item(13, 1).Value = 'mens long compression shorts'
$ws.Cells.Item(14, 1).Value = 'basketball knee sleeve'
$ws.Cells.Item(15, 1).Value = 'running pants for men'
$ws.Cells.Item(16, 1).Value = 'tall mens tights'
$ws.Cells.Item(17, 1).Value = 'baseball compression'
$ws.Cells.Item(18, 1).Value = 'protective pad'
$ws.Cells.Item(19, 1).Value = 'baseball youth pants'
$ws.Cells.Item(20, 1).Value = 'kneeling pad for exercise'
$ws.Cells.Item(21, 1).Value = 'basketball compression sleeve youth'
$ws.Cells.Item(22, 1).Value = 'volleyball pads for girls'
$ws.Cells.Item(23, 1).Value = 'sheer protectors'
$ws.Cells.Item(24, 1).Value = 'compression pants with pouch'
$ws.Cells.Item(25, 1).Value = 'protect knee'
$ws.Cells.Item(26, 1).Value = 'calf sleeve weight'
$ws.Cells.Item(27, 1).Value = 'mens cycling pants with padding'
$ws.Cells.Item(28, 1).Value = 'hip pads sports'
$ws.Cells.Item(29, 1).Value = 'baseball shorts'
$ws.Cells.Item(30, 1).Value = 'air knee pads'
$ws.Cells.Item(31, 1).Value = 'knee pads work construction'
$ws.Cells.Item(32, 1).Value = 'knee sleeves for basketball pair'
$ws.Cells.Item(33, 1).Value = 'compression tights men pack'
$ws.Cells.Item(34, 1).Value = 'football leggings for boys'
$ws.Cells.Item(35, 1).Value = 'construction pants for men'
$ws.Cells.Item(36, 1).Value = 'knee compression sleeve youth'
$ws.Cells.Item(37, 1).Value = 'knee pads thin'
$ws.Cells.Item(38, 1).Value = 'mens leggings pouch'
$ws.Cells.Item(39, 1).Value = 'running pants men tights'
$ws.Cells.Item(40, 1).Value = 'patella band'
$ws.Cells.Item(41, 1).Value = 'youth knee sleeve'
$ws.Cells.Item(42, 1).Value = 'calf compression sleeve spandex'
$ws.Cells.Item(43, 1).Value = 'youth baseball compression sleeve'
$ws.Cells.Item(44, 1).Value = 'protective pads'
$ws.Cells.Item(45, 1).Value = 'mens leggings cold'
$ws.Cells.Item(46, 1).Value = 'youth boys compression'
$ws.Cells.Item(47, 1).Value = 'little boys baseball pants size 5'
$ws.Cells.Item(48, 1).Value = 'mens running tights pants'
$ws.Cells.Item(49, 1).Value = 'men tights sports'
$ws.Cells.Item(50, 1).Value = 'snowboarding padded shorts'
$ws.Cells.Item(51, 1).Value = 'baseball dirt'
$ws.Cells.Item(52, 1).Value = 'boy compression pants'
$ws.Cells.Item(53, 1).Value = 'youth kneepads'
$ws.Cells.Item(54, 1).Value = 'best construction knee pads'
$ws.Cells.Item(55, 1).Value = 'fall cycling pants'
$ws.Cells.Item(56, 1).Value = 'knee pads mountain biking'
$ws.Cells.Item(57, 1).Value = 'football hip pads'
$ws.Cells.Item(58, 1).Value = 'catchers leg guards adult'
$ws.Cells.Item(59, 1).Value = 'flexible knee pad'
$ws.Cells.Item(60, 1).Value = 'basketball shorts pack of 5'
$ws.Cells.Item(61, 1).Value = 'volleyball mens shorts'
$ws.Cells.Item(62, 1).Value = 'knee pads girls'
$ws.Cells.Item(63, 1).Value = 'bees knees'
$ws.Cells.Item(64, 1).Value = 'sport leggings'
$ws.Cells.Item(65, 1).Value = 'girls lacrosse shorts'
$ws.Cells.Item(66, 1).Value = 'best knee pads'
$ws.Cells.Item(67, 1).Value = 'calf pads'
$ws.Cells.Item(68, 1).Value = 'burns baseball'
$ws.Cells.Item(69, 1).Value = 'boys workout leggings'
$ws.Cells.Item(70, 1).Value = '6 short pants'
$ws.Cells.Item(71, 1).Value = 'basketball leg sleeve youth'
$ws.Cells.Item(72, 1).Value = 'calf compression leggings'
$ws.Cells.Item(73, 1).Value = 'capri leggings with mesh'
$ws.Cells.Item(74, 1).Value = 'basketball sleeve for youth'
$ws.Cells.Item(75, 1).Value = 'extra thick knee pads'
$ws.Cells.Item(76, 1).Value = 'girdle football adult'
$ws.Cells.Item(77, 1).Value = 'indoor pants'
$ws.Cells.Item(78, 1).Value = 'youth girls knee pads'
$ws.Cells.Item(79, 1).Value = 'mens volleyball knee sleeve'
$ws.Cells.Item(80, 1).Value = 'knees for men'
$ws.Cells.Item(81, 1).Value = 'patella band knee'
$ws.Cells.Item(82, 1).Value = 'football knee sleeves'
$ws.Cells.Item(83, 1).Value = 'basketball sleeve for men'
$ws.Cells.Item(84, 1).Value = 'paintball pads and protection'
$ws.Cells.Item(85, 1).Value = 'sports leggings boys'
$ws.Cells.Item(86, 1).Value = 'football calf sleeve'
$ws.Cells.Item(87, 1).Value = '$5 and below'
$ws.Cells.Item(88, 1).Value = 'capri tights'
$ws.Cells.Item(89, 1).Value = 'compression bands for knees'
$ws.Cells.Item(90, 1).Value = 'mens compression pants black'
$ws.Cells.Item(91, 1).Value = 'playing ball on running water'
$ws.Cells.Item(92, 1).Value = 'youth large baseball pants'
$ws.Cells.Item(93, 1).Value = 'girls spandex shorts black volleyball'
$ws.Cells.Item(94, 1).Value = 'football short tights'
$ws.Cells.Item(95, 1).Value = 'thigh pads'
$ws.Cells.Item(96, 1).Value = 'tights and leggings'
$ws.Cells.Item(97, 1).Value = 'professional knee pad'
$ws.Cells.Item(98, 1).Value = 'short youth baseball pants'
$ws.Cells.Item(99, 1).Value = 'youth calf compression sleeve'
$ws.Cells.Item(100, 1).Value = 'knee guards'
